$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Mon Oct 09 23:50:50 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:51:08 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:51:26 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:51:44 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:52:02 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 09 23:52:20 EDT 2023"

$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Mon Oct 09 23:52:39 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:52:55 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:53:11 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:53:28 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:53:45 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 09 23:54:01 EDT 2023"
$ws.Range("B8").Value = "Mon Oct 09 23:54:18 EDT 2023"
$ws.Range("B9").Value = "Mon Oct 09 23:54:35 EDT 2023"
$ws.Range("B10").Value = "Mon Oct 09 23:54:51 EDT 2023"
$ws.Range("B11").Value = "Mon Oct 09 23:55:08 EDT 2023"
$ws.Range("B12").Value = "Mon Oct 09 23:55:24 EDT 2023"

$ws = $wb.Worksheets.Item("Extension")
$ws.Range("B2").Value = "Mon Oct 09 23:55:41 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:55:59 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:56:16 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:56:34 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:56:51 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 09 23:57:08 EDT 2023"

$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Mon Oct 09 23:57:26 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 23:57:43 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 23:58:00 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 23:58:17 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 09 23:58:35 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 09 23:58:53 EDT 2023"
$ws.Range("B8").Value = "Mon Oct 09 23:59:10 EDT 2023"
$ws.Range("B9").Value = "Mon Oct 09 23:59:27 EDT 2023"
$ws.Range("B10").Value = "Mon Oct 09 23:59:45 EDT 2023"
$ws.Range("B11").Value = "Tue Oct 10 00:00:03 EDT 2023"
$ws.Range("B12").Value = "Tue Oct 10 00:00:20 EDT 2023"
$ws.Range("B13").Value = "Tue Oct 10 00:00:38 EDT 2023"
$ws.Range("B14").Value = "Tue Oct 10 00:00:56 EDT 2023"
$ws.Range("B15").Value = "Tue Oct 10 00:01:14 EDT 2023"
$ws.Range("B16").Value = "Tue Oct 10 00:01:31 EDT 2023"

$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Tue Oct 10 00:02:10 EDT 2023"
$ws.Range("B3").Value = "Tue Oct 10 00:02:28 EDT 2023"
$ws.Range("B4").Value = "Tue Oct 10 00:02:46 EDT 2023"
$ws.Range("B5").Value = "Tue Oct 10 00:03:04 EDT 2023"
$ws.Range("B6").Value = "Tue Oct 10 00:03:23 EDT 2023"

$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Tue Oct 10 00:03:42 EDT 2023"
$ws.Range("B3").Value = "Tue Oct 10 00:04:08 EDT 2023"
$ws.Range("B4").Value = "Tue Oct 10 00:04:33 EDT 2023"
$ws.Range("B5").Value = "Tue Oct 10 00:04:59 EDT 2023"
$ws.Range("B6").Value = "Tue Oct 10 00:05:25 EDT 2023"

$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Tue Oct 10 00:01:49 EDT 2023"
